$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.206.33"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.642.07"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.254"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "1.872.92"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.636.96"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.542"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "27.206.51"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  +3.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "1.305.79"
$ws.Range("E34").Value = "  +3.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.551"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.857"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.25%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("E42").Value = "  +5.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").Value = "1.782.87"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
